$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValoresVenais")

# The "Tabela4" Excel Table (A1:AQ47) gains one more data row (2026),
# mirroring the previous row's formatting.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Carry the row-47 formatting down into the freshly added row 48
# (two groups, matching the "Ano/Decreto" columns vs. the currency columns).
$ws.Range("A47:B47").Copy()
$ws.Range("A48:B48").PasteSpecial(-4122) | Out-Null
$ws.Range("C47:AQ47").Copy()
$ws.Range("C48:AQ48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the 2026 figures. "Decreto" (B48) is intentionally left blank —
# no decree number published yet for this row.
$ws.Cells.Item(48,1).Value = 2026
$ws.Cells.Item(48,3).Value = 357.68
$ws.Cells.Item(48,4).Value = 286.17
$ws.Cells.Item(48,5).Value = 214.58
$ws.Cells.Item(48,6).Value = 171.54
$ws.Cells.Item(48,7).Value = 143.03
$ws.Cells.Item(48,8).Value = 114.49
$ws.Cells.Item(48,9).Value = 100
$ws.Cells.Item(48,10).Value = 85.72
$ws.Cells.Item(48,11).Value = 80.39
$ws.Cells.Item(48,12).Value = 0
$ws.Cells.Item(48,13).Value = 0
$ws.Cells.Item(48,14).Value = 0
$ws.Cells.Item(48,15).Value = 0
$ws.Cells.Item(48,16).Value = 0
$ws.Cells.Item(48,17).Value = 0
$ws.Cells.Item(48,18).Value = 0
$ws.Cells.Item(48,19).Value = 0
$ws.Cells.Item(48,20).Value = 0
$ws.Cells.Item(48,21).Value = 0
$ws.Cells.Item(48,22).Value = 0
$ws.Cells.Item(48,23).Value = 160.19
$ws.Cells.Item(48,24).Value = 356.01
$ws.Cells.Item(48,25).Value = 0
$ws.Cells.Item(48,26).Value = 640.79999999999995
$ws.Cells.Item(48,27).Value = 889.98
$ws.Cells.Item(48,28).Value = 1246.02
$ws.Cells.Item(48,29).Value = 0
$ws.Cells.Item(48,30).Value = 267.01
$ws.Cells.Item(48,31).Value = 0
$ws.Cells.Item(48,32).Value = 587.4
$ws.Cells.Item(48,33).Value = 818.78
$ws.Cells.Item(48,34).Value = 0
$ws.Cells.Item(48,35).Value = 0
$ws.Cells.Item(48,36).Value = 240.28
$ws.Cells.Item(48,37).Value = 0
$ws.Cells.Item(48,38).Value = 516.19000000000005
$ws.Cells.Item(48,39).Value = 676.38
$ws.Cells.Item(48,40).Value = 0
$ws.Cells.Item(48,41).Value = 0
$ws.Cells.Item(48,42).Value = 0
$ws.Cells.Item(48,43).Value = 3.31

# Move the view/selection down to the newly entered row, as the author did.
$ws.Activate() | Out-Null
$ws.Range("B48").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
